$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.413.42'
$ws.Range('E2').Value = '  -4.01%  '
$ws.Range('D3').Value = '2.906.01'
$ws.Range('E3').Value = '  -4.16%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''501.54'
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('D6').Value = '''132.28'
$ws.Range('E6').Value = '  -5.50%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '''0.419'
$ws.Range('E8').Value = '  -5.17%  '
$ws.Range('D9').Value = '''7.15'
$ws.Range('E9').Value = '  -4.58%  '
$ws.Range('E10').Value = '  -6.49%  '
$ws.Range('D11').Value = '''0.347'
$ws.Range('E11').Value = '  -5.66%  '
$ws.Range('D12').Value = '3.408.00'
$ws.Range('E12').Value = '  -4.00%  '
$ws.Range('E13').Value = '  -3.97%  '
$ws.Range('D14').Value = '''25.65'
$ws.Range('E14').Value = '  -4.26%  '
$ws.Range('D15').Value = '''0.0000158'
$ws.Range('E15').Value = '  -5.28%  '
$ws.Range('D16').Value = '55.403.48'
$ws.Range('E16').Value = '  -4.01%  '
$ws.Range('D17').Value = '''5.97'
$ws.Range('E17').Value = '  -4.32%  '
$ws.Range('D18').Value = '2.908.39'
$ws.Range('E18').Value = '  -3.96%  '
$ws.Range('D19').Value = '''12.47'
$ws.Range('E19').Value = '  -3.53%  '
$ws.Range('D20').Value = '''7.64'
$ws.Range('E20').Value = '  -4.88%  '
$ws.Range('D21').Value = '''312.63'
$ws.Range('E21').Value = '  -6.57%  '
$ws.Range('D22').Value = '''1.00'
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').Value = '''0.481'
$ws.Range('E23').Value = '  -3.95%  '
$ws.Range('D24').Value = '''62.48'
$ws.Range('E24').Value = '  -3.31%  '
$ws.Range('D25').Value = '3.029.84'
$ws.Range('E25').Value = '  -3.91%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = '''0.159'
$ws.Range('E27').Value = '  -5.84%  '
$ws.Range('D28').Value = '0.0₃0832'
$ws.Range('E28').Value = '  -10.45%  '
$ws.Range('D29').Value = '''6.36'
$ws.Range('E29').Value = '  -6.99%  '
$ws.Range('D30').Value = '''6.89'
$ws.Range('E30').Value = '  -8.08%  '
$ws.Range('D31').Value = '''1.76'
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('D32').Value = '''19.65'
$ws.Range('E32').Value = '  -5.80%  '
$ws.Range('D33').Value = '''1.12'
$ws.Range('E33').Value = '  -8.08%  '
$ws.Range('D34').Value = '''149.31'
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('D35').Value = '''4.35'
$ws.Range('E35').Value = '  -8.00%  '
$ws.Range('D36').Value = '''5.55'
$ws.Range('E36').Value = '  -5.58%  '
$ws.Range('D37').Value = '''24.29'
$ws.Range('E37').Value = '  -2.66%  '
$ws.Range('E38').Value = '  -8.25%  '
$ws.Range('D39').Value = '''0.0642'
$ws.Range('E39').Value = '  -6.53%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').Value = '''36.23'
$ws.Range('E41').Value = '  -3.21%  '
$ws.Range('D42').Value = '''3.68'
$ws.Range('E42').Value = '  -5.38%  '
$ws.Range('D43').Value = '''0.634'
$ws.Range('E43').Value = '  -3.36%  '
$ws.Range('D44').Value = '2.089.10'
$ws.Range('E44').Value = '  -9.78%  '
$ws.Range('D45').Value = '''1.32'
$ws.Range('E45').Value = '  -7.62%  '
$ws.Range('D46').Value = '''5.87'
$ws.Range('E46').Value = '  -2.76%  '
$ws.Range('D47').Value = '''0.905'
$ws.Range('E47').Value = '  -8.75%  '
$ws.Range('E48').Value = '  -3.88%  '
$ws.Range('D49').Value = '''18.51'
$ws.Range('E49').Value = '  -6.59%  '
$ws.Range('D50').Value = '''0.0833'
$ws.Range('E50').Value = '  -7.20%  '
$ws.Range('D51').Value = '''1.66'
$ws.Range('E51').Value = '  -9.67%  '
